# Applies the "Changes in the tests, error strategy" commit:
#  - Rename sheets: "TypesOK" -> "Data OK", "TypesWithErrors" -> "Data With Errors"
#  - "Data OK" (sheet 1): bold the header row; it's no longer the active tab/
#    selection becomes the whole used range (A1:J6) instead of I7.
#  - "Data With Errors" (sheet 2): gains the same 10-column layout as
#    "Data OK" (Int/Decimal/Float/Nullable.../String/Date/Nullable Date/Bool),
#    with row 3 mostly blanked out and row 5 replaced by literal error values
#    ("a" in most columns, "A" in the Bool column). It becomes the active tab,
#    with D6 selected.

$wb = $excel.ActiveWorkbook

$wsOk = $wb.Worksheets.Item(1)
$wsErr = $wb.Worksheets.Item(2)

# --- Rename sheets ---
$wsOk.Name = "Data OK"
$wsErr.Name = "Data With Errors"

# --- "Data OK": bold header row ---
$wsOk.Range("A1:J1").Font.Bold = $true

# --- "Data With Errors": rebuild to mirror "Data OK"'s layout ---
$wsErr.Cells.Clear()
$wsOk.Range("A1:J6").Copy($wsErr.Range("A1"))

# Row 3 only keeps the three formatted-but-empty cells (D3, H3, I3)
$wsErr.Rows.Item(3).ClearContents()

# Row 5 becomes the "error" row: bool column gets "A", everything else "a"
$wsErr.Range("J5").Value = "A"
$wsErr.Range("A5:I5").Value = "a"

# Column widths now need to fit the wider (10-column) layout, same as "Data OK".
# (internal stored width = ColumnWidth + 5/6)
$wsErr.Columns.Item(3).ColumnWidth = 12.7109375 - 0.8333333333333333
$wsErr.Columns.Item(4).ColumnWidth = 18.85546875 - 0.8333333333333333
$wsErr.Columns.Item(5).ColumnWidth = 23.7109375 - 0.8333333333333333
$wsErr.Columns.Item(6).ColumnWidth = 20.85546875 - 0.8333333333333333
$wsErr.Columns.Item(7).ColumnWidth = 13.42578125 - 0.8333333333333333
$wsErr.Columns.Item(8).ColumnWidth = 12.42578125 - 0.8333333333333333
$wsErr.Columns.Item(9).ColumnWidth = 20.5703125 - 0.8333333333333333

# --- Selection / active tab state ---
[void]$wsOk.Select()
[void]$wsOk.Range("A1:J6").Select()

[void]$wsErr.Select()
[void]$wsErr.Range("D6").Select()
